$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: add new column S with the 2022 header ---
# Copy R4's format (existing style used for year headers) onto the new S4 cell
# so that the new cell reuses the existing style index instead of creating a
# duplicate style, then set its value.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S4").Value = 2022

# --- Row 5: update Q5/R5 values and add new S5 value ---
# Q5 switches to the percent-number-format style used by R5, with an updated value.
$ws.Range("R5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Q5").Value = 91.892815141492093

# R5 keeps its existing style, only the value changes.
$ws.Range("R5").Value = 101.53074848578628

# S5 is a brand-new cell; give it the same format as R5, then set its value.
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S5").Value = 109.27053140096621

$excel.CutCopyMode = 0

# --- Update the saved selection to match the author's final cursor position ---
$null = $ws.Range("T5").Select()
